$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 2 (mirrors the author manually selecting the row before deleting it)
$ws.Range("A2:XFD2").Select()

# Delete the entire row 2 (the "CP1-3513SJCT-ND" 3.5mm audio jack line that
# Dylan Thorner noted couldn't be found as a footprint / wasn't in the
# schematic). Everything below shifts up by one row.
$ws.Rows("2").Delete()
